$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '67.746.90'
$ws.Range("E2").Value = '  +0.54%  '

# Row 3
$ws.Range("D3").Value = '2.495.15'
$ws.Range("E3").Value = '  -2.28%  '

# Row 4
$ws.Range("E4").Value = '  -0.02%  '

# Row 5
$ws.Range("D5").Value = '''591.86'

# Row 6
$ws.Range("D6").Value = '''173.84'
$ws.Range("E6").Value = '  +0.43%  '

# Row 7
$ws.Range("E7").Value = '  -0.06%  '

# Row 8
$ws.Range("E8").Value = '  -1.07%  '

# Row 9
$ws.Range("D9").Value = '2.494.07'
$ws.Range("E9").Value = '  -2.32%  '

# Row 10
$ws.Range("E10").Value = '  -0.19%  '

# Row 11
$ws.Range("E11").Value = '  +1.79%  '

# Row 12
$ws.Range("D12").Value = '''5.09'
$ws.Range("E12").Value = '  -1.55%  '

# Row 13
$ws.Range("E13").Value = '  -2.89%  '

# Row 14
$ws.Range("D14").Value = '''26.28'
$ws.Range("E14").Value = '  -3.17%  '

# Row 15
$ws.Range("D15").Value = '2.949.35'
$ws.Range("E15").Value = '  -1.85%  '

# Row 16
$ws.Range("D16").Value = '''0.0000177'
$ws.Range("E16").Value = '  -1.44%  '

# Row 17
$ws.Range("D17").Value = '67.553.46'
$ws.Range("E17").Value = '  +0.55%  '

# Row 18
$ws.Range("D18").Value = '2.499.46'
$ws.Range("E18").Value = '  -2.08%  '

# Row 19
$ws.Range("B19").Value = 'Uniswap'
$ws.Range("C19").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D19").Value = '''8.04'
$ws.Range("E19").Value = '  +1.13%  '

# Row 20
$ws.Range("B20").Value = 'Chainlink'
$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D20").Value = '''11.73'
$ws.Range("E20").Value = '  +2.39%  '

# Row 21
$ws.Range("D21").Value = '''364.93'
$ws.Range("E21").Value = '  +2.18%  '

# Row 22
$ws.Range("D22").Value = '''4.13'
$ws.Range("E22").Value = '  -2.49%  '

# Row 23
$ws.Range("D23").Value = '''4.56'
$ws.Range("E23").Value = '  -2.94%  '

# Row 24
$ws.Range("B24").Value = 'Dai'
$ws.Range("C24").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D24").Value = '''1.00'
$ws.Range("E24").Value = '  -0.07%  '

# Row 25
$ws.Range("B25").Value = 'Litecoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D25").Value = '''71.17'
$ws.Range("E25").Value = '  +1.34%  '

# Row 26
$ws.Range("E26").Value = '  -6.13%  '

# Row 27
$ws.Range("D27").Value = '''9.85'
$ws.Range("E27").Value = '  -2.27%  '

# Row 28
$ws.Range("D28").Value = '''1.00'
$ws.Range("E28").Value = '  -0.21%  '

# Row 29
$ws.Range("D29").Value = '2.610.39'

# Row 30
$ws.Range("D30").Value = '0.0₃0966'
$ws.Range("E30").Value = '  -3.52%  '

# Row 31
$ws.Range("D31").Value = '''532.56'
$ws.Range("E31").Value = '  -0.62%  '

# Row 32
$ws.Range("E32").Value = '  -0.36%  '

# Row 33
$ws.Range("D33").Value = '''1.87'
$ws.Range("E33").Value = '  -0.06%  '

# Row 34
$ws.Range("E34").Value = '  -4.77%  '

# Row 35
$ws.Range("D35").Value = '''1.00'
$ws.Range("E35").Value = '  -0.05%  '

# Row 36
$ws.Range("E36").Value = '  -4.70%  '

# Row 37
$ws.Range("D37").Value = '''158.82'
$ws.Range("E37").Value = '  +0.92%  '

# Row 38
$ws.Range("D38").Value = '''1.43'
$ws.Range("E38").Value = '  -3.93%  '

# Row 39
$ws.Range("D39").Value = '''18.61'
$ws.Range("E39").Value = '  -1.22%  '

# Row 40
$ws.Range("E40").Value = '  +0.90%  '

# Row 41
$ws.Range("E41").Value = '  -1.64%  '

# Row 42
$ws.Range("B42").Value = 'PolygonEcosystemToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D42").Value = '''0.349'
$ws.Range("E42").Value = '  -2.66%  '

# Row 43
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D43").Value = '''5.11'
$ws.Range("E43").Value = '  -2.14%  '

# Row 44
$ws.Range("D44").Value = '''0.998'
$ws.Range("E44").Value = '  -0.27%  '

# Row 45
$ws.Range("D45").Value = '''2.50'
$ws.Range("E45").Value = '  -2.47%  '

# Row 46
$ws.Range("D46").Value = '''144.92'
$ws.Range("E46").Value = '  -4.32%  '

# Row 47
$ws.Range("E47").Value = '  -1.53%  '

# Row 48
$ws.Range("D48").Value = '''0.548'
$ws.Range("E48").Value = '  -3.47%  '

# Row 49
$ws.Range("D49").Value = '0.0₆0272'
$ws.Range("E49").Value = '  -4.03%  '

# Row 50
$ws.Range("D50").Value = '''1.70'
$ws.Range("E50").Value = '  -1.81%  '

# Row 51
$ws.Range("E51").Value = '  -1.97%  '
